# ============================================================
# Edit: add "2022-Q3" quarterly fund-holding data
#  1. Summary sheet "总计": insert a new row (row 2) with the
#     2022-Q3 aggregate (count=40, holding value=4.64), shifting
#     the existing quarters down and renumbering the index column.
#  2. Insert a brand-new worksheet named "2022-Q3" (placed right
#     before "2022-Q2") containing the per-fund holding detail.
# ============================================================

$wb = $excel.ActiveWorkbook

# ---------- 1. Summary sheet ("总计") ----------
$summary = $wb.Worksheets.Item("总计")

# Make room for the new quarter at row 2 (existing rows 2-8 -> 3-9).
$summary.Rows.Item(2).Insert()

# Re-fetch a fresh reference (sheet collection shifted by Insert) and
# copy the row formatting down from the row that now holds the old
# "2022-Q2" record so the new row matches the table's look.
$summary = $wb.Worksheets.Item("总计")
$summary.Range("A3:D3").Copy($summary.Range("A2:D2"))

# Write the new 2022-Q3 summary values.
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 40
$summary.Cells.Item(2, 4).Value = 4.64

# The first column is a running 0-based row index - renumber rows 2..9.
for ($i = 0; $i -lt 8; $i++) {
    $summary.Cells.Item(2 + $i, 1).Value = $i
}

# ---------- 2. New worksheet "2022-Q3" ----------
$template = $wb.Worksheets.Item("2022-Q2")
$added = $wb.Worksheets.Add($template)
$added.Name = "2022-Q3"

# Re-fetch both sheets by name - a stale reference after Add() does not
# reliably carry formatting operations in this host.
$q3 = $wb.Worksheets.Item("2022-Q3")
$template = $wb.Worksheets.Item("2022-Q2")

# Header row (B1:H1) - copy text + formatting from the template sheet.
$template.Range("B1:H1").Copy($q3.Range("B1:H1"))

$q3data = @(
    @(0, '010178', '大成企业能力驱动混合A', '30.94', '81.17', '3.30', '1.0210', 8),
    @(1, '011230', '创金合信数字经济主题股票C', '15.89', '92.35', '5.42', '0.8612', 1),
    @(2, '011229', '创金合信数字经济主题股票A', '15.42', '92.35', '5.42', '0.8358', 1),
    @(3, '009447', '财通资管科技创新一年定期开放混合', '8.69', '94.05', '6.10', '0.5301', 3),
    @(4, '090004', '大成精选增值混合', '9.37', '81.16', '2.98', '0.2792', 9),
    @(5, '020009', '国泰金鹏蓝筹混合', '5.16', '89.77', '3.32', '0.1713', 2),
    @(6, '161631', '融通中证人工智能主题指数（LOF）A', '4.54', '92.94', '2.53', '0.1149', 10),
    @(7, '165523', '信诚中证信息安全指数（LOF）A', '3.21', '93.83', '2.83', '0.0908', 9),
    @(8, '515400', '富国中证大数据产业ETF', '2.11', '99.08', '4.02', '0.0848', 8),
    @(9, '001144', '大成互联网思维混合', '2.35', '86.55', '3.28', '0.0771', 9),
    @(10, '516510', '易方达中证云计算与大数据主题ETF', '1.70', '98.74', '3.79', '0.0644', 8),
    @(11, '161628', '融通中证云计算与大数据主题指数（LOF）A', '1.52', '94.31', '3.78', '0.0575', 8),
    @(12, '560660', '新华中证云计算50ETF', '1.20', '97.03', '4.36', '0.0523', 7),
    @(13, '007853', '华商计算机行业量化股票', '1.67', '91.78', '2.65', '0.0443', 10),
    @(14, '009907', '湘财长泽灵活配置混合A', '0.93', '87.32', '4.40', '0.0409', 9),
    @(15, '159613', '嘉实中证信息安全主题ETF', '1.24', '99.34', '2.98', '0.0370', 9),
    @(16, '159738', '华泰柏瑞中证沪港深云计算产业ETF', '0.71', '95.55', '3.45', '0.0245', 9),
    @(17, '516630', '华夏中证云计算与大数据主题ETF', '0.64', '99.46', '3.82', '0.0244', 8),
    @(18, '159739', '鹏华中证云计算与大数据主题ETF', '0.62', '97.89', '3.76', '0.0233', 8),
    @(19, '014130', '融通中证云计算与大数据主题指数（LOF）C', '0.58', '94.31', '3.78', '0.0219', 8),
    @(20, '516000', '华夏中证大数据产业ETF', '0.50', '97.60', '3.98', '0.0199', 8),
    @(21, '010179', '大成企业能力驱动混合C', '0.59', '81.17', '3.30', '0.0195', 8),
    @(22, '517390', '天弘中证沪港深云计算产业ETF', '0.52', '99.26', '3.48', '0.0181', 9),
    @(23, '009908', '湘财长泽灵活配置混合C', '0.39', '87.32', '4.40', '0.0172', 9),
    @(24, '009239', '融通中证人工智能主题指数（LOF）C', '0.67', '92.94', '2.53', '0.0170', 10),
    @(25, '168701', '合煦智远国证香蜜湖金融科技指数（LOF）A', '0.47', '92.07', '3.37', '0.0158', 6),
    @(26, '012432', '国投瑞银安泰混合C', '1.00', '32.06', '1.51', '0.0151', 7),
    @(27, '159890', '招商中证云计算与大数据主题ETF', '0.36', '98.52', '3.76', '0.0135', 8),
    @(28, '516700', '华宝中证大数据产业ETF', '0.23', '95.99', '3.91', '0.0090', 8),
    @(29, '012019', '国投瑞银安泽混合A', '0.49', '32.69', '1.63', '0.0080', 6),
    @(30, '015201', '创金合信动态平衡混合C', '0.23', '65.33', '3.11', '0.0072', 5),
    @(31, '013083', '信诚中证信息安全指数（LOF）C', '0.25', '93.83', '2.83', '0.0071', 9),
    @(32, '015200', '创金合信动态平衡混合A', '0.16', '65.33', '3.11', '0.0050', 5),
    @(33, '168702', '合煦智远国证香蜜湖金融科技指数（LOF）C', '0.13', '92.07', '3.37', '0.0044', 6),
    @(34, '014543', '汇添富中证沪港深云计算产业指数A', '0.09', '93.00', '3.28', '0.0030', 9),
    @(35, '011494', '华泰紫金丰和偏债混合发起A', '0.12', '31.18', '1.63', '0.0020', 1),
    @(36, '012020', '国投瑞银安泽混合C', '0.10', '32.69', '1.63', '0.0016', 6),
    @(37, '014544', '汇添富中证沪港深云计算产业指数C', '0.04', '93.00', '3.28', '0.0013', 9),
    @(38, '011495', '华泰紫金丰和偏债混合发起C', '0.04', '31.18', '1.63', '0.0007', 1),
    @(39, '012431', '国投瑞银安泰混合A', '0.00', '32.06', '1.51', '0', 7)
)

for ($i = 0; $i -lt $q3data.Count; $i++) {
    $r = 2 + $i
    $row = $q3data[$i]

    # Column A (index) - numeric, styled like the template's index column.
    $template.Range("A2").Copy($q3.Cells.Item($r, 1))
    $q3.Cells.Item($r, 1).Value = $row[0]

    # Columns B-G are stored as text in the source data (leading zeros in
    # fund codes must survive) - force text with a leading apostrophe.
    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = "'" + $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]

    # G41 (last row) is the sole numeric exception in the source data.
    if ($i -eq ($q3data.Count - 1)) {
        $q3.Cells.Item($r, 7).Value = 0
    } else {
        $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    }

    # Column H (rank) - numeric.
    $q3.Cells.Item($r, 8).Value = $row[7]
}

Write-Host "2022-Q3 sheet populated with" $q3data.Count "rows"
